# edit.ps1 -- Word COM-interop script that applies the changes described
# in the target diff to CIS4930Report.docx.
#
# Summary of changes:
#  1. Author block is rewritten: the single "Miles Brosz, Natalee Sama, ..."
#     line is split out so each author gets their own line, a 4th author
#     (Aidan Mahoney) is added, and the affiliation / course info / "Prof.
#     Mithila" / date block is re-ordered underneath.
#  2. "( Miami" -> "(Miami" (stray space after the opening parenthesis).
#  3. A new closing sentence is appended to the paragraph describing the
#     learning-model choice ("We chose this style of learning model
#     because we found the Linear-Regression idea ... easiest to
#     implement.")

$d = $word.ActiveDocument

function Set-RunFormatting($rng) {
    # Match the document's body text formatting (12pt / sz=24, szCs=24)
    $rng.Font.Size = 12
    $rng.Font.SizeBi = 12
}

# ---------------------------------------------------------------------
# 1. Author / title block
# ---------------------------------------------------------------------

# "Miles Brosz, Natalee Sama, ..." -> "Miles Brosz,"
$d.Content.Find.Execute("Brosz, Natalee Sama, ...", $true, $false, $false, $false, $false, $true, 1, $false, "Brosz,", 2) | Out-Null

# "Florida State University" -> " Natalee Sama,"  (becomes the 2nd author line)
$d.Content.Find.Execute("Florida State University", $true, $false, $false, $false, $false, $true, 1, $false, " Natalee Sama,", 2) | Out-Null

# "CIS4930; Special Topics" -> " Tyler Zuluaga"  (becomes the 4th author line)
$d.Content.Find.Execute("CIS4930; Special Topics", $true, $false, $false, $false, $false, $true, 1, $false, " Tyler Zuluaga", 2) | Out-Null

# "Mithila" -> "Florida State University"  (affiliation line moves down)
$d.Content.Find.Execute("Mithila", $true, $false, $false, $false, $false, $true, 1, $false, "Florida State University", 2) | Out-Null

# Clear the old "4/11/25" run -- the paragraph stays but becomes empty; the
# date is re-added a few paragraphs later on.
$d.Content.Find.Execute("4/11/25", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# The empty paragraph right after "Florida State University" (which now
# reads " Natalee Sama,") is paragraph 10 -- insert the 3rd author,
# "Aidan Mahoney,", there.
$p10 = $d.Paragraphs(10)
$p10.Range.InsertAfter("Aidan Mahoney,")
Set-RunFormatting $p10.Range

# Paragraph 13 is now the empty paragraph left behind by clearing
# "4/11/25". Insert three new paragraphs (course info, instructor, date)
# right after it, before the "Title2" paragraph that follows.
$p13 = $d.Paragraphs(13)

$p13.Range.InsertParagraphAfter()
$pCourse = $d.Paragraphs(14)
$pCourse.Range.InsertAfter("CIS4930; Special Topics")
Set-RunFormatting $pCourse.Range

$pCourse.Range.InsertParagraphAfter()
$pProf = $d.Paragraphs(15)
$pProf.Range.InsertAfter("Prof. ")
Set-RunFormatting $pProf.Range
$pProf.Range.InsertAfter("Mithila")
Set-RunFormatting $pProf.Range

$pProf.Range.InsertParagraphAfter()
$pDate = $d.Paragraphs(16)
$pDate.Range.InsertAfter("4/11/25")
Set-RunFormatting $pDate.Range

# ---------------------------------------------------------------------
# 2. "( Miami" -> "(Miami"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("( Miami", $true, $false, $false, $false, $false, $true, 1, $false, "(Miami", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Append the new closing sentence about why this learning model was
#    chosen, right after "...from the regression training."
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("from the regression training. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $targetPara = $rng.Paragraphs(1).Range
    $targetPara.InsertAfter("We chose this style of learning model because we found the Linear-Regression idea that it is based from to be the most intuitive and easiest to implement.")
    Set-RunFormatting $targetPara
}

Write-Host "Done applying edits."
